$d = $word.ActiveDocument

# --- Edit 1: Bryan paragraph - fix missing space "Ooijende" -> "Ooijen de" ---
$d.Content.Find.Execute("Ooijende klassen", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Ooijen de klassen", 2)

# --- Edit 2: Chanan paragraph - append new status update after "Chanan:" ---
$rng = $d.Content
$rng.Find.Execute("Chanan:", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)
$vtab = [char]11
$tabchar = [char]9

$rng.InsertAfter($vtab + 'Ik heb samen met Bryan Baan de klassen Equalize, Histogram, LineDetection, MedianFilter en Threshold (versie van Yusuf')
$rng.Collapse(0)
$rng.InsertAfter(' herschreven) ')
$rng.Collapse(0)
$rng.InsertAfter('gemaakt. Vrijdag 07/03/2014 heb ik daar aan gewerkt van 11:00 – 17:00 en zaterdag van 10:00 – 16:00 en 23:00 tot zondagochtend 06:00.')
$rng.Collapse(0)
$rng.InsertAfter(' ')
$rng.Collapse(0)
$rng.InsertAfter('“')
$rng.Collapse(0)
$rng.InsertAfter('De te testen algoritmen – Algoritme ')
$rng.Collapse(0)
$rng.InsertAfter('1')
$rng.Collapse(0)
$rng.InsertAfter('”')
$rng.Collapse(0)
$rng.InsertAfter(' ')
$rng.Collapse(0)
$rng.InsertAfter('en ')
$rng.Collapse(0)
$rng.InsertAfter('“')
$rng.Collapse(0)
$rng.InsertAfter('W')
$rng.Collapse(0)
$rng.InsertAfter('at gaan we testen')
$rng.Collapse(0)
$rng.InsertAfter('”')
$rng.Collapse(0)

# zero-width "_GoBack" bookmark right after the closing quote
$markStart = $rng.Start
$rng.InsertAfter('Q')
$bmRange = $d.Range($markStart, $rng.End)
$d.Bookmarks.Add("_GoBack", $bmRange)
$bmRange.Text = ""
$rng = $d.Range($bmRange.End, $bmRange.End)

$rng.InsertAfter(' heb ik geschreven.')
$rng.Collapse(0)
$rng.InsertAfter($tabchar)
$rng.Collapse(0)
